# SGDS-RFC03.docx: "Eliminando firma de la solicitud de cambio 03."
#
# The document ends with:
#   <w:p/>                                   (empty)
#   <w:p> "Saavedra Monterrey Max Bruno/ Jefe de proyecto" </w:p>   (signature name/title)
#   <w:p> <drawing of the signature image> </w:p>                  (signature image)
#   <w:p/>                                   (empty)
#   <w:sectPr> ... </w:sectPr>
#
# The signature (both the name/title text and the scanned-signature picture)
# must be removed, leaving three empty paragraphs before the section
# properties - i.e. the two signature paragraphs collapse into one empty
# paragraph, which then sits between the two pre-existing empty paragraphs.

$d = $word.ActiveDocument

# --- 1. Clear the "Saavedra Monterrey Max Bruno/ Jefe de proyecto" text ---
# Locate the paragraph by its distinctive text, expand the range to the
# whole paragraph (this includes the trailing paragraph mark), then shrink
# back by one character so the paragraph mark itself is preserved and only
# the run contents are removed.
$sigRange = $d.Content
$found = $sigRange.Find.Execute("Saavedra Monterrey Max Bruno", $true, $false, $false, $false,
                                 $false, $true, 1, $false, "", 0)
if ($found) {
    $sigRange.Expand(4) | Out-Null   # wdParagraph - grow to the full paragraph incl. mark
    $sigParaStart = $sigRange.Start
    $sigRange.End = $sigRange.End - 1  # keep the paragraph mark itself
    $sigRange.Text = ""
} else {
    $sigParaStart = -1
}

# --- 2. Delete the signature picture (the drawing anchored in the next paragraph) ---
for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $shp = $d.Shapes.Item($i)
    if ($shp.Name -like "Imagen*") {
        $shp.Delete()
    }
}

# --- 3. Merge the now-empty signature paragraph with the (already empty)
#        paragraph that used to hold the picture, so only one empty
#        paragraph remains where the two signature paragraphs used to be. ---
if ($sigParaStart -ge 0) {
    $probe = $d.Range($sigParaStart, $sigParaStart + 1)
    $sigPara = $probe.Paragraphs.Item(1)
    $markRange = $d.Range($sigPara.Range.End - 1, $sigPara.Range.End)
    $markRange.Delete()
}
